# REVER_DailyTracker_BALRAJ.xlsx -- "Add files via upload"
#
# The NOV-2020 sheet gets 4 new task rows (10-13) describing RPA RLOGIC
# work, plus several more blank-but-formatted rows (14-18) below them
# before the trailing legend block (rows 19-23). Column D is widened to
# fit the new (longer) text, and the sheet's scroll/selection moves down
# to where the new data was typed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NOV-2020")
$ws.Activate()

# --- 1. Pre-format the new block (A10:G18) ------------------------------
# Rows 10-16 already exist as empty, custom-formatted rows; rows 17-18 are
# brand new. Give the whole block the same bordered/general-number-format
# look the rest of the table uses (the style carried by G2, a plain
# bordered cell), then layer the Date format (carried by B2) onto B10 and
# the Percent format (carried by E2) onto E10:E13 -- matching exactly
# which of the new cells end up holding real values.
$ws.Range("G2").Copy() | Out-Null
$ws.Range("A10:G18").PasteSpecial(-4122) | Out-Null

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B10").PasteSpecial(-4122) | Out-Null

$ws.Range("E2").Copy() | Out-Null
$ws.Range("E10:E13").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- 2. Fill in the new task rows ---------------------------------------
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 44140
$ws.Range("C10").Value = "RPA RLOGIC "
$ws.Range("D10").Value = "1. Log has been implemented at DRS download"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = "Completed"

$ws.Range("D11").Value = "2. Converted Python script to Exe(Trial Version) for the DRS download and tested and running smoothly"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = "Completed"

$ws.Range("D12").Value = "3. Converted Python script to Exe (Trial Version) for RPA Management, tested and running smoothly"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = "Completed"

$ws.Range("D13").Value = "4. Sending email with attachments has been done for RPA Management template"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = "Completed"

# --- 3. Widen column D so the longer comments fit -----------------------
$ws.Columns.Item(4).ColumnWidth = 85.5

# --- 4. Move the view/selection down to the newly-entered data ----------
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("F13").Select() | Out-Null
